$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (N_Calib_1=20, N_Calib_2=40)
$ws.Range("C2").Value = 0.2540379456918681
$ws.Range("D2").Value = 0.8018254903282214

# Row 3 (N_Calib_1=20, N_Calib_2=60)
$ws.Range("C3").Value = 1.040128729193758
$ws.Range("D3").Value = 0.3095768363330005

# Row 4 (N_Calib_1=20, N_Calib_2=100)
$ws.Range("C4").Value = 0.1390219754357667
$ws.Range("D4").Value = 0.8906975569796585

# Row 5 (N_Calib_1=20, N_Calib_2=200)
$ws.Range("C5").Value = 1.9232197405741
$ws.Range("D5").Value = 0.06749042111471493
$ws.Range("G5").Value = "No"

# Row 6 (N_Calib_1=40, N_Calib_2=60)
$ws.Range("C6").Value = 0.7785164688858091
$ws.Range("D6").Value = 0.4445579049958635

# Row 7 (N_Calib_1=40, N_Calib_2=100)
$ws.Range("C7").Value = -0.2182268958908977
$ws.Range("D7").Value = 0.8292638929800735

# Row 8 (N_Calib_1=40, N_Calib_2=200)
$ws.Range("C8").Value = 1.841268366022929
$ws.Range("D8").Value = 0.07910751346193567

# Row 9 (N_Calib_1=60, N_Calib_2=100)
$ws.Range("C9").Value = -1.206823610606248
$ws.Range("D9").Value = 0.2403168368763708

# Row 10 (N_Calib_1=60, N_Calib_2=200)
$ws.Range("C10").Value = 0.6715237013300431
$ws.Range("D10").Value = 0.5088775258138791

# Row 11 (N_Calib_1=100, N_Calib_2=200)
$ws.Range("C11").Value = 1.764001878372826
$ws.Range("D11").Value = 0.09161502812602706
$ws.Range("G11").Value = "No"
